# "Add files via upload" — rename the course title on slides 1-6 from
# "13 jQuery Server App" to "13 jQuery Server App: Part 1".
$p = $ppt.ActivePresentation

for ($i = 1; $i -le 6; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Text = "13 jQuery Server App: Part 1"
}
